$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new timetable rows for course SC2207 (rows 99-148)
$ws.Range("A99").Value = "SC2207"
$ws.Range("B99").Value = "10321"
$ws.Range("C99").Value = "TUT"
$ws.Range("D99").Value = "THU"
$ws.Range("E99").Value = "10:30"
$ws.Range("F99").Value = "11:20"
$ws.Range("G99").Value = "SCSA"
$ws.Range("A100").Value = "SC2207"
$ws.Range("B100").Value = "10321"
$ws.Range("C100").Value = "LAB"
$ws.Range("D100").Value = "MON"
$ws.Range("E100").Value = "10:30"
$ws.Range("F100").Value = "12:20"
$ws.Range("G100").Value = "SCSA"
$ws.Range("H100").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A101").Value = "SC2207"
$ws.Range("B101").Value = "10322"
$ws.Range("C101").Value = "TUT"
$ws.Range("D101").Value = "TUE"
$ws.Range("E101").Value = "10:30"
$ws.Range("F101").Value = "11:20"
$ws.Range("G101").Value = "SCSB"
$ws.Range("A102").Value = "SC2207"
$ws.Range("B102").Value = "10322"
$ws.Range("C102").Value = "LAB"
$ws.Range("D102").Value = "MON"
$ws.Range("E102").Value = "12:30"
$ws.Range("F102").Value = "14:20"
$ws.Range("G102").Value = "SCSB"
$ws.Range("G102").NumberFormat = "h:mm"
$ws.Range("H102").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Rows(102).RowHeight = 15
$ws.Range("A103").Value = "SC2207"
$ws.Range("B103").Value = "10323"
$ws.Range("C103").Value = "TUT"
$ws.Range("D103").Value = "THU"
$ws.Range("E103").Value = "10:30"
$ws.Range("F103").Value = "11:20"
$ws.Range("G103").Value = "SCSC"
$ws.Range("A104").Value = "SC2207"
$ws.Range("B104").Value = "10323"
$ws.Range("C104").Value = "LAB"
$ws.Range("D104").Value = "MON"
$ws.Range("E104").Value = "10:30"
$ws.Range("F104").Value = "12:20"
$ws.Range("G104").Value = "SCSC"
$ws.Range("H104").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A105").Value = "SC2207"
$ws.Range("B105").Value = "10324"
$ws.Range("C105").Value = "TUT"
$ws.Range("D105").Value = "TUE"
$ws.Range("E105").Value = "15:30"
$ws.Range("F105").Value = "16:20"
$ws.Range("G105").Value = "SCSD"
$ws.Range("A106").Value = "SC2207"
$ws.Range("B106").Value = "10324"
$ws.Range("C106").Value = "LAB"
$ws.Range("D106").Value = "WED"
$ws.Range("E106").Value = "12:30"
$ws.Range("F106").Value = "14:20"
$ws.Range("G106").Value = "SCSD"
$ws.Range("H106").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A107").Value = "SC2207"
$ws.Range("B107").Value = "10325"
$ws.Range("C107").Value = "TUT"
$ws.Range("D107").Value = "TUE"
$ws.Range("E107").Value = "15:30"
$ws.Range("F107").Value = "16:20"
$ws.Range("G107").Value = "SCSE"
$ws.Range("A108").Value = "SC2207"
$ws.Range("B108").Value = "10325"
$ws.Range("C108").Value = "LAB"
$ws.Range("D108").Value = "WED"
$ws.Range("E108").Value = "12:30"
$ws.Range("F108").Value = "14:20"
$ws.Range("G108").Value = "SCSE"
$ws.Range("H108").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A109").Value = "SC2207"
$ws.Range("B109").Value = "10326"
$ws.Range("C109").Value = "TUT"
$ws.Range("D109").Value = "THU"
$ws.Range("E109").Value = "12:30"
$ws.Range("F109").Value = "13:20"
$ws.Range("G109").Value = "SCSF"
$ws.Range("A110").Value = "SC2207"
$ws.Range("B110").Value = "10326"
$ws.Range("C110").Value = "LAB"
$ws.Range("D110").Value = "THU"
$ws.Range("E110").Value = "08:30"
$ws.Range("F110").Value = "10:20"
$ws.Range("G110").Value = "SCSF"
$ws.Range("H110").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A111").Value = "SC2207"
$ws.Range("B111").Value = "10327"
$ws.Range("C111").Value = "TUT"
$ws.Range("D111").Value = "MON"
$ws.Range("E111").Value = "09:30"
$ws.Range("F111").Value = "10:20"
$ws.Range("G111").Value = "SCSX"
$ws.Range("A112").Value = "SC2207"
$ws.Range("B112").Value = "10327"
$ws.Range("C112").Value = "LAB"
$ws.Range("D112").Value = "MON"
$ws.Range("E112").Value = "12:30"
$ws.Range("F112").Value = "14:20"
$ws.Range("G112").Value = "SCSX"
$ws.Range("H112").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A113").Value = "SC2207"
$ws.Range("B113").Value = "10328"
$ws.Range("C113").Value = "TUT"
$ws.Range("D113").Value = "THU"
$ws.Range("E113").Value = "09:30"
$ws.Range("F113").Value = "10:20"
$ws.Range("G113").Value = "SCSY"
$ws.Range("A114").Value = "SC2207"
$ws.Range("B114").Value = "10328"
$ws.Range("C114").Value = "LAB"
$ws.Range("D114").Value = "THU"
$ws.Range("E114").Value = "12:30"
$ws.Range("F114").Value = "14:20"
$ws.Range("G114").Value = "SCSY"
$ws.Range("H114").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A115").Value = "SC2207"
$ws.Range("B115").Value = "10329"
$ws.Range("C115").Value = "TUT"
$ws.Range("D115").Value = "TUE"
$ws.Range("E115").Value = "15:30"
$ws.Range("F115").Value = "16:20"
$ws.Range("G115").Value = "SCS1"
$ws.Range("A116").Value = "SC2207"
$ws.Range("B116").Value = "10329"
$ws.Range("C116").Value = "LAB"
$ws.Range("D116").Value = "FRI"
$ws.Range("E116").Value = "13:30"
$ws.Range("F116").Value = "15:20"
$ws.Range("G116").Value = "SCS1"
$ws.Range("H116").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A117").Value = "SC2207"
$ws.Range("B117").Value = "10330"
$ws.Range("C117").Value = "TUT"
$ws.Range("D117").Value = "TUE"
$ws.Range("E117").Value = "15:30"
$ws.Range("F117").Value = "16:20"
$ws.Range("G117").Value = "SCS2"
$ws.Range("A118").Value = "SC2207"
$ws.Range("B118").Value = "10330"
$ws.Range("C118").Value = "LAB"
$ws.Range("D118").Value = "FRI"
$ws.Range("E118").Value = "13:30"
$ws.Range("F118").Value = "15:20"
$ws.Range("G118").Value = "SCS2"
$ws.Range("H118").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A119").Value = "SC2207"
$ws.Range("B119").Value = "10331"
$ws.Range("C119").Value = "TUT"
$ws.Range("D119").Value = "THU"
$ws.Range("E119").Value = "13:30"
$ws.Range("F119").Value = "14:20"
$ws.Range("G119").Value = "SCS3"
$ws.Range("A120").Value = "SC2207"
$ws.Range("B120").Value = "10331"
$ws.Range("C120").Value = "LAB"
$ws.Range("D120").Value = "THU"
$ws.Range("E120").Value = "10:30"
$ws.Range("F120").Value = "12:20"
$ws.Range("G120").Value = "SCS3"
$ws.Range("H120").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A121").Value = "SC2207"
$ws.Range("B121").Value = "10332"
$ws.Range("C121").Value = "TUT"
$ws.Range("D121").Value = "THU"
$ws.Range("E121").Value = "13:30"
$ws.Range("F121").Value = "14:20"
$ws.Range("G121").Value = "SCS4"
$ws.Range("A122").Value = "SC2207"
$ws.Range("B122").Value = "10332"
$ws.Range("C122").Value = "LAB"
$ws.Range("D122").Value = "THU"
$ws.Range("E122").Value = "10:30"
$ws.Range("F122").Value = "12:20"
$ws.Range("G122").Value = "SCS4"
$ws.Range("H122").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A123").Value = "SC2207"
$ws.Range("B123").Value = "10333"
$ws.Range("C123").Value = "TUT"
$ws.Range("D123").Value = "WED"
$ws.Range("E123").Value = "13:30"
$ws.Range("F123").Value = "14:20"
$ws.Range("G123").Value = "SCS5"
$ws.Range("A124").Value = "SC2207"
$ws.Range("B124").Value = "10333"
$ws.Range("C124").Value = "LAB"
$ws.Range("D124").Value = "WED"
$ws.Range("E124").Value = "10:30"
$ws.Range("F124").Value = "12:20"
$ws.Range("G124").Value = "SCS5"
$ws.Range("H124").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A125").Value = "SC2207"
$ws.Range("B125").Value = "10334"
$ws.Range("C125").Value = "TUT"
$ws.Range("D125").Value = "WED"
$ws.Range("E125").Value = "13:30"
$ws.Range("F125").Value = "14:20"
$ws.Range("G125").Value = "SCS6"
$ws.Range("A126").Value = "SC2207"
$ws.Range("B126").Value = "10334"
$ws.Range("C126").Value = "LAB"
$ws.Range("D126").Value = "WED"
$ws.Range("E126").Value = "10:30"
$ws.Range("F126").Value = "12:20"
$ws.Range("G126").Value = "SCS6"
$ws.Range("H126").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A127").Value = "SC2207"
$ws.Range("B127").Value = "10335"
$ws.Range("C127").Value = "TUT"
$ws.Range("D127").Value = "TUE"
$ws.Range("E127").Value = "13:30"
$ws.Range("F127").Value = "14:20"
$ws.Range("G127").Value = "SCMB"
$ws.Range("A128").Value = "SC2207"
$ws.Range("B128").Value = "10335"
$ws.Range("C128").Value = "LAB"
$ws.Range("D128").Value = "TUE"
$ws.Range("E128").Value = "14:30"
$ws.Range("F128").Value = "16:20"
$ws.Range("G128").Value = "SCMB"
$ws.Range("H128").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A129").Value = "SC2207"
$ws.Range("B129").Value = "10336"
$ws.Range("C129").Value = "TUT"
$ws.Range("D129").Value = "TUE"
$ws.Range("E129").Value = "13:30"
$ws.Range("F129").Value = "14:20"
$ws.Range("G129").Value = "SCMC"
$ws.Range("A130").Value = "SC2207"
$ws.Range("B130").Value = "10336"
$ws.Range("C130").Value = "LAB"
$ws.Range("D130").Value = "TUE"
$ws.Range("E130").Value = "14:30"
$ws.Range("F130").Value = "16:20"
$ws.Range("G130").Value = "SCMC"
$ws.Range("H130").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A131").Value = "SC2207"
$ws.Range("B131").Value = "10337"
$ws.Range("C131").Value = "TUT"
$ws.Range("D131").Value = "THU"
$ws.Range("E131").Value = "09:30"
$ws.Range("F131").Value = "10:20"
$ws.Range("G131").Value = "STA1"
$ws.Range("A132").Value = "SC2207"
$ws.Range("B132").Value = "10337"
$ws.Range("C132").Value = "LAB"
$ws.Range("D132").Value = "THU"
$ws.Range("E132").Value = "10:30"
$ws.Range("F132").Value = "12:20"
$ws.Range("G132").Value = "STA1"
$ws.Range("H132").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A133").Value = "SC2207"
$ws.Range("B133").Value = "10338"
$ws.Range("C133").Value = "TUT"
$ws.Range("D133").Value = "WED"
$ws.Range("E133").Value = "12:30"
$ws.Range("F133").Value = "13:20"
$ws.Range("G133").Value = "ACDA"
$ws.Range("A134").Value = "SC2207"
$ws.Range("B134").Value = "10338"
$ws.Range("C134").Value = "LAB"
$ws.Range("D134").Value = "THU"
$ws.Range("E134").Value = "12:30"
$ws.Range("F134").Value = "14:20"
$ws.Range("G134").Value = "ACDA"
$ws.Range("H134").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A135").Value = "SC2207"
$ws.Range("B135").Value = "10339"
$ws.Range("C135").Value = "TUT"
$ws.Range("D135").Value = "WED"
$ws.Range("E135").Value = "12:30"
$ws.Range("F135").Value = "13:20"
$ws.Range("G135").Value = "BACF1"
$ws.Range("A136").Value = "SC2207"
$ws.Range("B136").Value = "10339"
$ws.Range("C136").Value = "LAB"
$ws.Range("D136").Value = "TUE"
$ws.Range("E136").Value = "12:30"
$ws.Range("F136").Value = "14:20"
$ws.Range("G136").Value = "BACF1"
$ws.Range("H136").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A137").Value = "SC2207"
$ws.Range("B137").Value = "10340"
$ws.Range("C137").Value = "TUT"
$ws.Range("D137").Value = "WED"
$ws.Range("E137").Value = "13:30"
$ws.Range("F137").Value = "14:20"
$ws.Range("G137").Value = "BACF2"
$ws.Range("A138").Value = "SC2207"
$ws.Range("B138").Value = "10340"
$ws.Range("C138").Value = "LAB"
$ws.Range("D138").Value = "TUE"
$ws.Range("E138").Value = "12:30"
$ws.Range("F138").Value = "14:20"
$ws.Range("G138").Value = "BACF2"
$ws.Range("H138").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A139").Value = "SC2207"
$ws.Range("B139").Value = "10341"
$ws.Range("C139").Value = "TUT"
$ws.Range("D139").Value = "FRI"
$ws.Range("E139").Value = "09:30"
$ws.Range("F139").Value = "10:20"
$ws.Range("G139").Value = "ECDS1"
$ws.Range("A140").Value = "SC2207"
$ws.Range("B140").Value = "10341"
$ws.Range("C140").Value = "LAB"
$ws.Range("D140").Value = "WED"
$ws.Range("E140").Value = "10:30"
$ws.Range("F140").Value = "12:20"
$ws.Range("G140").Value = "ECDS1"
$ws.Range("H140").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A141").Value = "SC2207"
$ws.Range("B141").Value = "10342"
$ws.Range("C141").Value = "TUT"
$ws.Range("D141").Value = "THU"
$ws.Range("E141").Value = "11:30"
$ws.Range("F141").Value = "12:20"
$ws.Range("G141").Value = "ECDS2"
$ws.Range("A142").Value = "SC2207"
$ws.Range("B142").Value = "10342"
$ws.Range("C142").Value = "LAB"
$ws.Range("D142").Value = "WED"
$ws.Range("E142").Value = "10:30"
$ws.Range("F142").Value = "12:20"
$ws.Range("G142").Value = "ECDS2"
$ws.Range("H142").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A143").Value = "SC2207"
$ws.Range("B143").Value = "10343"
$ws.Range("C143").Value = "TUT"
$ws.Range("D143").Value = "WED"
$ws.Range("E143").Value = "12:30"
$ws.Range("F143").Value = "13:20"
$ws.Range("G143").Value = "MACS"
$ws.Range("A144").Value = "SC2207"
$ws.Range("B144").Value = "10343"
$ws.Range("C144").Value = "LAB"
$ws.Range("D144").Value = "TUE"
$ws.Range("E144").Value = "09:30"
$ws.Range("F144").Value = "11:20"
$ws.Range("G144").Value = "MACS"
$ws.Range("H144").Value = "Teaching Wk2,4,6,8,10,12"
$ws.Range("A145").Value = "SC2207"
$ws.Range("B145").Value = "10344"
$ws.Range("C145").Value = "TUT"
$ws.Range("D145").Value = "FRI"
$ws.Range("E145").Value = "15:30"
$ws.Range("F145").Value = "16:20"
$ws.Range("G145").Value = "REP"
$ws.Range("A146").Value = "SC2207"
$ws.Range("B146").Value = "10344"
$ws.Range("C146").Value = "LAB"
$ws.Range("D146").Value = "MON"
$ws.Range("E146").Value = "10:30"
$ws.Range("F146").Value = "12:20"
$ws.Range("G146").Value = "REP"
$ws.Range("H146").Value = "Teaching Wk1,3,5,7,9,11,13"
$ws.Range("A147").Value = "SC2207"
$ws.Range("B147").Value = "10345"
$ws.Range("C147").Value = "TUT"
$ws.Range("D147").Value = "TUE"
$ws.Range("E147").Value = "10:30"
$ws.Range("F147").Value = "11:20"
$ws.Range("G147").Value = "SCSG"
$ws.Range("A148").Value = "SC2207"
$ws.Range("B148").Value = "10345"
$ws.Range("C148").Value = "LAB"
$ws.Range("D148").Value = "MON"
$ws.Range("E148").Value = "12:30"
$ws.Range("F148").Value = "14:20"
$ws.Range("G148").Value = "SCSG"
$ws.Range("H148").Value = "Teaching Wk1,3,5,7,9,11,13"

# Update the view: scroll position and selection to match the new data range
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 88
$win.ScrollColumn = 1
$ws.Range("A103:H148").Select() | Out-Null
